# cryptos.xlsx refresh -- GitHub Actions scrape update
#
# Rewrites the "Price" (D) and "Volume(1h)" (E) columns for every coin row
# (rows 2-51) with freshly scraped values, and also fixes the ranking swap
# between FraxShare and TheSandbox (rows 42-43: Coin/Link/Price/Volume all
# move together).
#
# Price-column quirk: every price in this sheet is stored as literal TEXT,
# never a number (e.g. "1.002", "0.06336", "30.487.16"). Assigning a
# numeric-looking string straight to Range.Value lets Excel's type
# inference coerce it into a real number, so for any value that looks like
# a plain decimal we prepend a "'" (apostrophe) the same way a user typing
# into the grid would force text entry. Values Excel could never parse as
# a number anyway (e.g. "30.487.16", with two dots) are assigned as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    if ($Text -match '^-?\d+(\.\d+)?$') {
        # Looks like a plain number to Excel -- force text entry, the way
        # typing a leading apostrophe in the grid does.
        $Range.Value = "'" + $Text
    } else {
        $Range.Value = $Text
    }
}

Set-TextValue $ws.Range('D2') '30.487.16'
$ws.Range('E2').Value = '  +0.40%  '

Set-TextValue $ws.Range('D3') '1.854.00'
$ws.Range('E3').Value = '  -0.39%  '

Set-TextValue $ws.Range('D4') '1.002'
$ws.Range('E4').Value = '  +0.33%  '

Set-TextValue $ws.Range('D5') '233.24'
$ws.Range('E5').Value = '  -0.77%  '

Set-TextValue $ws.Range('D6') '1.002'
$ws.Range('E6').Value = '  +0.27%  '

Set-TextValue $ws.Range('D7') '0.4696'
$ws.Range('E7').Value = '  -0.74%  '

Set-TextValue $ws.Range('D8') '0.2738'
$ws.Range('E8').Value = '  -0.52%  '

Set-TextValue $ws.Range('D9') '0.06336'
$ws.Range('E9').Value = '  -1.65%  '

Set-TextValue $ws.Range('D10') '17.41'
$ws.Range('E10').Value = '  +6.36%  '

Set-TextValue $ws.Range('D11') '1.874.34'
$ws.Range('E11').Value = '  +1.32%  '

Set-TextValue $ws.Range('D12') '0.07431'
$ws.Range('E12').Value = '  +0.12%  '

Set-TextValue $ws.Range('D13') '5.098'
$ws.Range('E13').Value = '  +1.72%  '

Set-TextValue $ws.Range('D14') '84.53'
$ws.Range('E14').Value = '  -1.38%  '

Set-TextValue $ws.Range('D15') '0.6256'
$ws.Range('E15').Value = '  -1.61%  '

Set-TextValue $ws.Range('D16') '30.506.26'
$ws.Range('E16').Value = '  +0.63%  '

Set-TextValue $ws.Range('D17') '241.66'
$ws.Range('E17').Value = '  +3.60%  '

$ws.Range('E18').Value = '  +0.22%  '

Set-TextValue $ws.Range('D19') '12.71'
$ws.Range('E19').Value = '  -0.58%  '

Set-TextValue $ws.Range('D20') '0.000007325'
$ws.Range('E20').Value = '  -1.01%  '

Set-TextValue $ws.Range('D21') '1.003'
$ws.Range('E21').Value = '  +0.45%  '

Set-TextValue $ws.Range('D22') '4.976'
$ws.Range('E22').Value = '  -0.96%  '

Set-TextValue $ws.Range('D23') '6.002'
$ws.Range('E23').Value = '  -0.35%  '

Set-TextValue $ws.Range('D24') '9.280'
$ws.Range('E24').Value = '  -0.24%  '

Set-TextValue $ws.Range('D25') '164.23'
$ws.Range('E25').Value = '  -0.98%  '

Set-TextValue $ws.Range('D26') '18.05'
$ws.Range('E26').Value = '  +0.47%  '

Set-TextValue $ws.Range('D27') '1.887'
$ws.Range('E27').Value = '  -0.42%  '

Set-TextValue $ws.Range('D28') '0.1016'
$ws.Range('E28').Value = '  -2.41%  '

Set-TextValue $ws.Range('D29') '1.381'
$ws.Range('E29').Value = '  -0.35%  '

Set-TextValue $ws.Range('D30') '4.042'
$ws.Range('E30').Value = '  -2.48%  '

Set-TextValue $ws.Range('D31') '3.850'
$ws.Range('E31').Value = '  -2.44%  '

Set-TextValue $ws.Range('D32') '0.04904'
$ws.Range('E32').Value = '  -0.21%  '

Set-TextValue $ws.Range('D33') '1.145'
$ws.Range('E33').Value = '  -0.50%  '

Set-TextValue $ws.Range('D34') '0.7051'
$ws.Range('E34').Value = '  -3.40%  '

Set-TextValue $ws.Range('D35') '2.710'
$ws.Range('E35').Value = '  +0.63%  '

Set-TextValue $ws.Range('D36') '0.01901'
$ws.Range('E36').Value = '  -0.41%  '

Set-TextValue $ws.Range('D37') '2.682'
$ws.Range('E37').Value = '  +1.34%  '

Set-TextValue $ws.Range('D38') '0.8784'
$ws.Range('E38').Value = '  -3.40%  '

Set-TextValue $ws.Range('D39') '1.978'
$ws.Range('E39').Value = '  -0.16%  '

Set-TextValue $ws.Range('D40') '105.12'
$ws.Range('E40').Value = '  -0.76%  '

Set-TextValue $ws.Range('D41') '1.002'
$ws.Range('E41').Value = '  +0.21%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D42') '5.531'
$ws.Range('E42').Value = '  -0.21%  '

$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D43') '0.4068'
$ws.Range('E43').Value = '  -1.38%  '

Set-TextValue $ws.Range('D44') '7.235'
$ws.Range('E44').Value = '  +0.93%  '

Set-TextValue $ws.Range('D45') '62.96'
$ws.Range('E45').Value = '  +2.67%  '

Set-TextValue $ws.Range('D46') '0.1202'
$ws.Range('E46').Value = '  -0.69%  '

Set-TextValue $ws.Range('D47') '8.632'
$ws.Range('E47').Value = '  -1.41%  '

Set-TextValue $ws.Range('D48') '33.42'
$ws.Range('E48').Value = '  +0.06%  '

Set-TextValue $ws.Range('D49') '0.05542'
$ws.Range('E49').Value = '  -0.89%  '

Set-TextValue $ws.Range('D50') '1.358'
$ws.Range('E50').Value = '  -3.71%  '

Set-TextValue $ws.Range('D51') '0.3672'
$ws.Range('E51').Value = '  -1.40%  '
